# "Generate Report for Archive"
#  - Update status text "Ready for handoff" -> "In Translation" on all three
#    sheets (Overview, zh-cn, de-de).
#  - Narrow the "Latest Handoff Datetime" / status-date columns (Overview!E:F,
#    zh-cn!C, de-de!C) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

# ---- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("C3").Value = "In Translation"

$wsZh.Columns.Item(3).ColumnWidth = 13.4101845877511

# ---- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("C3").Value = "In Translation"

$wsDe.Columns.Item(3).ColumnWidth = 13.4101845877511
